$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.574712643678161
$ws.Range("C2").Value = 0.016602809706258
$ws.Range("D2").Value = 0.0600255427841635
$ws.Range("E2").Value = 0.945083014048531
$ws.Range("F2").Value = 0.0114942528735632
$ws.Range("G2").Value = 0.943805874840358
$ws.Range("H2").Value = 0.0140485312899106
$ws.Range("I2").Value = 0.731800766283525
$ws.Range("J2").Value = 0.0408684546615581
$ws.Range("K2").Value = 0.0395913154533844
$ws.Range("L2").Value = 0.0357598978288633
$ws.Range("M2").Value = 0.840357598978289
$ws.Range("N2").Value = 0.0102171136653895
$ws.Range("P2").Value = 0.033205619412516
$ws.Range("Q2").Value = 0.9272030651341
$ws.Range("R2").Value = 0.00638569604086846
$ws.Range("S2").Value = 0.00383141762452107
$ws.Range("U2").Value = 0.00638569604086846
$ws.Range("W2").Value = 0.0114942528735632
$ws.Range("X2").Value = 0.00766283524904215
$ws.Range("B3").Value = 0.0664112388250319
$ws.Range("C3").Value = 0.348659003831418
$ws.Range("D3").Value = 0.842911877394636
$ws.Range("E3").Value = 0.0446998722860792
$ws.Range("F3").Value = 0.00127713920817369
$ws.Range("H3").Value = 0.0689655172413793
$ws.Range("I3").Value = 0.0357598978288633
$ws.Range("J3").Value = 0.130268199233716
$ws.Range("K3").Value = 0.920817369093231
$ws.Range("L3").Value = 0.960408684546616
$ws.Range("M3").Value = 0.109833971902937
$ws.Range("N3").Value = 0.919540229885057
$ws.Range("O3").Value = 0.0153256704980843
$ws.Range("P3").Value = 0.00127713920817369
$ws.Range("R3").Value = 0.851851851851852
$ws.Range("S3").Value = 0.983397190293742
$ws.Range("T3").Value = 0.053639846743295
$ws.Range("U3").Value = 0.00766283524904215
$ws.Range("V3").Value = 0.0217113665389527
$ws.Range("W3").Value = 0.033205619412516
$ws.Range("X3").Value = 0.00127713920817369
$ws.Range("B4").Value = 0.352490421455939
$ws.Range("C4").Value = 0.0255427841634738
$ws.Range("D4").Value = 0.00383141762452107
$ws.Range("E4").Value = 0.00383141762452107
$ws.Range("F4").Value = 0.945083014048531
$ws.Range("G4").Value = 0.0523627075351213
$ws.Range("H4").Value = 0.00255427841634738
$ws.Range("I4").Value = 0.0191570881226054
$ws.Range("J4").Value = 0.053639846743295
$ws.Range("K4").Value = 0.037037037037037
$ws.Range("L4").Value = 0.00127713920817369
$ws.Range("M4").Value = 0.00638569604086846
$ws.Range("N4").Value = 0.00383141762452107
$ws.Range("P4").Value = 0.962962962962963
$ws.Range("Q4").Value = 0.00383141762452107
$ws.Range("R4").Value = 0.134099616858238
$ws.Range("U4").Value = 0.00255427841634738
$ws.Range("V4").Value = 0.0523627075351213
$ws.Range("W4").Value = 0.952745849297573
$ws.Range("X4").Value = 0.954022988505747
$ws.Range("B5").Value = 0.00638569604086846
$ws.Range("C5").Value = 0.609195402298851
$ws.Range("D5").Value = 0.0906768837803321
$ws.Range("E5").Value = 0.00510855683269476
$ws.Range("F5").Value = 0.0421455938697318
$ws.Range("G5").Value = 0.00383141762452107
$ws.Range("H5").Value = 0.914431673052363
$ws.Range("I5").Value = 0.213282247765006
$ws.Range("J5").Value = 0.77522349936143
$ws.Range("K5").Value = 0.00127713920817369
$ws.Range("L5").Value = 0.00255427841634738
$ws.Range("M5").Value = 0.0434227330779055
$ws.Range("N5").Value = 0.0664112388250319
$ws.Range("O5").Value = 0.984674329501916
$ws.Range("P5").Value = 0.00255427841634738
$ws.Range("Q5").Value = 0.0689655172413793
$ws.Range("R5").Value = 0.00766283524904215
$ws.Range("S5").Value = 0.0127713920817369
$ws.Range("T5").Value = 0.946360153256705
$ws.Range("U5").Value = 0.983397190293742
$ws.Range("V5").Value = 0.925925925925926
$ws.Range("W5").Value = 0.00255427841634738
$ws.Range("X5").Value = 0.037037037037037
